$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.92828380041271
$ws.Range("C2").Value = 11.99753558411496
$ws.Range("E2").Value = 16.61415612428156
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 21.41482617422153
$ws.Range("H2").Value = 11.82415776307961
$ws.Range("I2").Value = 15.76941180630444
$ws.Range("B3").Value = 14.09274380400131
$ws.Range("C3").Value = 11.27412348347575
$ws.Range("E3").Value = 15.66065143845027
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 21.20914396109521
$ws.Range("H3").Value = 11.89768919172315
$ws.Range("I3").Value = 15.95849712796945
$ws.Range("B4").Value = 13.55434689630764
$ws.Range("C4").Value = 10.8038133412604
$ws.Range("E4").Value = 15.05010367973166
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 21.10357408479674
$ws.Range("H4").Value = 11.94770839264732
$ws.Range("I4").Value = 16.08217755868612
$ws.Range("B5").Value = 13.32874625092809
$ws.Range("C5").Value = 10.60561920418915
$ws.Range("E5").Value = 14.79525727178504
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 21.06578984053047
$ws.Range("H5").Value = 11.96930441872635
$ws.Range("I5").Value = 16.13447065011653
$ws.Range("B6").Value = 13.29091701990578
$ws.Range("C6").Value = 10.57231545182811
$ws.Range("E6").Value = 14.75258406190304
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 21.05983238409871
$ws.Range("H6").Value = 11.97296332283882
$ws.Range("I6").Value = 16.14326777136863
$ws.Range("B7").Value = 13.55132920276886
$ws.Range("C7").Value = 10.8011668507541
$ws.Range("E7").Value = 15.04669080934304
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 21.10304329375008
$ws.Range("H7").Value = 11.94799474957776
$ws.Range("I7").Value = 16.08287515798285
$ws.Range("B8").Value = 14.64555268363097
$ws.Range("C8").Value = 11.75354275716681
$ws.Range("E8").Value = 16.29073957883498
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 21.33962913575249
$ws.Range("H8").Value = 11.84849331251114
$ws.Range("I8").Value = 15.83302626674619
$ws.Range("B9").Value = 16.58379255138533
$ws.Range("C9").Value = 13.41291302460753
$ws.Range("E9").Value = 18.64827316954648
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 21.96570745839082
$ws.Range("H9").Value = 11.69255198805081
$ws.Range("I9").Value = 15.40388959855343
$ws.Range("B10").Value = 17.87493243604384
$ws.Range("C10").Value = 14.50450098711207
$ws.Range("E10").Value = 20.31743651428219
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 22.52001250217444
$ws.Range("H10").Value = 11.60259779595784
$ws.Range("I10").Value = 15.1266395604282
$ws.Range("B11").Value = 18.43255036554023
$ws.Range("C11").Value = 14.97339197676254
$ws.Range("E11").Value = 21.03458920740287
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 22.79132675816599
$ws.Range("H11").Value = 11.56717610430719
$ws.Range("I11").Value = 15.00899687625171
$ws.Range("B12").Value = 18.63937735329088
$ws.Range("C12").Value = 15.1469729129337
$ws.Range("E12").Value = 21.30013390598946
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 22.89668973989112
$ws.Range("H12").Value = 11.554566281735
$ws.Range("I12").Value = 14.96568846710669
$ws.Range("B13").Value = 18.59502682446886
$ws.Range("C13").Value = 15.10976600132007
$ws.Range("E13").Value = 21.24321150362037
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 22.87388354343485
$ws.Range("H13").Value = 11.55724608860642
$ws.Range("I13").Value = 14.97496019849636
$ws.Range("B14").Value = 18.44965321523435
$ws.Range("C14").Value = 14.98775228376908
$ws.Range("E14").Value = 21.05655627849213
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 22.79994328107629
$ws.Range("H14").Value = 11.56612250796307
$ws.Range("I14").Value = 15.00540888544453
$ws.Range("B15").Value = 18.3600422511718
$ws.Range("C15").Value = 14.9124974394373
$ws.Range("E15").Value = 20.94144096039223
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 22.75499009082952
$ws.Range("H15").Value = 11.57166460371734
$ws.Range("I15").Value = 15.0242217529819
$ws.Range("B16").Value = 17.8378874889913
$ws.Range("C16").Value = 14.47330101203381
$ws.Range("E16").Value = 20.26972372216343
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 22.50265589901021
$ws.Range("H16").Value = 11.60502444156808
$ws.Range("I16").Value = 15.1345001852798
$ws.Range("B17").Value = 17.50990409381095
$ws.Range("C17").Value = 14.19677965641909
$ws.Range("E17").Value = 19.84687551069397
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 22.35267315034438
$ws.Range("H17").Value = 11.62690721647409
$ws.Range("I17").Value = 15.20433884361781
$ws.Range("B18").Value = 17.31846094261863
$ws.Range("C18").Value = 14.03512633843731
$ws.Range("E18").Value = 19.59969409381714
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 22.26821906330711
$ws.Range("H18").Value = 11.64001043134008
$ws.Range("I18").Value = 15.24530502528894
$ws.Range("B19").Value = 17.25316333717803
$ws.Range("C19").Value = 13.97994579003831
$ws.Range("E19").Value = 19.51531968454354
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 22.239939433059
$ws.Range("H19").Value = 11.6445353253822
$ws.Range("I19").Value = 15.25931170741041
$ws.Range("B20").Value = 17.54510828712449
$ws.Range("C20").Value = 14.22648548618738
$ws.Range("E20").Value = 19.89229906811818
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 22.36845245877806
$ws.Range("H20").Value = 11.62452417536599
$ws.Range("I20").Value = 15.19682178237847
$ws.Range("B21").Value = 18.49247088483097
$ws.Range("C21").Value = 15.02369859572449
$ws.Range("E21").Value = 21.1115446794511
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 22.82159127559113
$ws.Range("H21").Value = 11.56349337694084
$ws.Range("I21").Value = 14.99643153341801
$ws.Range("B22").Value = 19.08636680468321
$ws.Range("C22").Value = 15.52154041087393
$ws.Range("E22").Value = 21.87328813007909
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 23.1329589946284
$ws.Range("H22").Value = 11.52829611178482
$ws.Range("I22").Value = 14.87270500679297
$ws.Range("B23").Value = 18.77171981728474
$ws.Range("C23").Value = 15.25795256534736
$ws.Range("E23").Value = 21.46993121613931
$ws.Range("F23").Value = 21.82633154475857
$ws.Range("G23").Value = 22.96543051484955
$ws.Range("H23").Value = 11.54664821229947
$ws.Range("I23").Value = 14.93807049593224
$ws.Range("B24").Value = 17.5292014410394
$ws.Range("C24").Value = 14.21306382068349
$ws.Range("E24").Value = 19.87177579057829
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 22.361313100824
$ws.Range("H24").Value = 11.62559992157031
$ws.Range("I24").Value = 15.20021770855998
$ws.Range("B25").Value = 16.08241519618651
$ws.Range("C25").Value = 12.98639045598958
$ws.Range("E25").Value = 17.99595921573317
$ws.Range("F25").Value = 18.34778573295691
$ws.Range("G25").Value = 21.77937139203071
$ws.Range("H25").Value = 11.73046816574957
$ws.Range("I25").Value = 8.896963201771335
